$wb = $excel.ActiveWorkbook

# --- Update "Conversión del día" note on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.78 = 10305.56 pesos`n✅ 10305.56 pesos = 2.8 = 939.3 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 360
$wsTasas.Range("O10").Value = 3710
$wsTasas.Range("N12").Value = 3675.48
$wsTasas.Range("O12").Value = 335
